$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 222.6
$ws.Range("I9").Value = 62.153847
$ws.Range("K9").Value = 62.153847
$ws.Range("M9").Value = 106.846153
$ws.Range("H33").Value = 28087
$ws.Range("I33").Value = 32222.688
$ws.Range("K33").Value = 32222.688
$ws.Range("M33").Value = -31993.688
$ws.Range("H132").Value = 3702.0386
$ws.Range("I132").Value = 3489
$ws.Range("J132").Value = 4873.75
$ws.Range("K132").Value = 10467
$ws.Range("L132").Value = 14621.25
$ws.Range("M132").Value = -7937
$ws.Range("N132").Value = -19681.25
$ws.Range("H133").Value = 89997
$ws.Range("J133").Value = 89997
$ws.Range("L133").Value = 89997
$ws.Range("N133").Value = -100117
$ws.Range("H138").Value = 2901.2
$ws.Range("J138").Value = 4215.6523
$ws.Range("L138").Value = 12646.9569
$ws.Range("N138").Value = -22926.9569

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1149.6786
$ws.Range("I2").Value = 1366.2941
$ws.Range("J2").Value = 814.9091
$ws.Range("K2").Value = 1366.2941
$ws.Range("L2").Value = 814.9091
$ws.Range("M2").Value = -1253.2941
$ws.Range("N2").Value = -1040.9091
$ws.Range("H62").Value = 54999
$ws.Range("J62").Value = 54999
$ws.Range("L62").Value = 54999
$ws.Range("N62").Value = -56247
$ws.Range("H65").Value = 54999
$ws.Range("J65").Value = 54999
$ws.Range("L65").Value = 164997
$ws.Range("N65").Value = -171237
$ws.Range("H74").Value = 3053.7646
$ws.Range("I74").Value = 1443.0741
$ws.Range("K74").Value = 1443.0741
$ws.Range("M74").Value = -569.0741
$ws.Range("H77").Value = 3053.7646
$ws.Range("I77").Value = 1443.0741
$ws.Range("K77").Value = 7215.3705
$ws.Range("M77").Value = -2847.3705
$ws.Range("H97").Value = 326.04166
$ws.Range("I97").Value = 301.13635
$ws.Range("K97").Value = 301.13635
$ws.Range("M97").Value = 194.86365
$ws.Range("H110").Value = 1187.1428
$ws.Range("I110").Value = 1187.1428
$ws.Range("K110").Value = 1187.1428
$ws.Range("M110").Value = 857.8571999999999
$ws.Range("H116").Value = 1149.6786
$ws.Range("I116").Value = 1366.2941
$ws.Range("J116").Value = 814.9091
$ws.Range("K116").Value = 1366.2941
$ws.Range("L116").Value = 814.9091
$ws.Range("M116").Value = 927.7058999999999
$ws.Range("N116").Value = -5402.9091
$ws.Range("H122").Value = 2034.1111
$ws.Range("I122").Value = 2034.1111
$ws.Range("K122").Value = 6102.3333
$ws.Range("M122").Value = -3652.3333
$ws.Range("H132").Value = 4359.533
$ws.Range("I132").Value = 1999.3636
$ws.Range("J132").Value = 10850
$ws.Range("K132").Value = 5998.0908
$ws.Range("L132").Value = 32550
$ws.Range("M132").Value = -3468.0908
$ws.Range("N132").Value = -37610

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1149.6786
$ws.Range("I3").Value = 1366.2941
$ws.Range("J3").Value = 814.9091
$ws.Range("K3").Value = 1366.2941
$ws.Range("L3").Value = 814.9091
$ws.Range("M3").Value = -1252.2941
$ws.Range("N3").Value = -1042.9091
$ws.Range("H99").Value = 2735
$ws.Range("I99").Value = 2544.5454
$ws.Range("K99").Value = 2544.5454
$ws.Range("M99").Value = -1046.5454
$ws.Range("H134").Value = 3977.0715
$ws.Range("I134").Value = 1469.3
$ws.Range("K134").Value = 4407.9
$ws.Range("M134").Value = -1872.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2478.6177
$ws.Range("I99").Value = 1595.3462
$ws.Range("J99").Value = 5349.25
$ws.Range("K99").Value = 1595.3462
$ws.Range("L99").Value = 5349.25
$ws.Range("M99").Value = -97.34619999999995
$ws.Range("N99").Value = -8345.25
$ws.Range("H105").Value = 2319.6
$ws.Range("I105").Value = 2313.7144
$ws.Range("J105").Value = 2333.3333
$ws.Range("K105").Value = 2313.7144
$ws.Range("L105").Value = 2333.3333
$ws.Range("M105").Value = -566.7143999999998
$ws.Range("N105").Value = -5827.3333
$ws.Range("H122").Value = 2069.625
$ws.Range("I122").Value = 1491.6666
$ws.Range("J122").Value = 3803.5
$ws.Range("K122").Value = 4474.9998
$ws.Range("L122").Value = 11410.5
$ws.Range("M122").Value = -2024.9998
$ws.Range("N122").Value = -16310.5
$ws.Range("H126").Value = 2478.6177
$ws.Range("I126").Value = 1595.3462
$ws.Range("J126").Value = 5349.25
$ws.Range("K126").Value = 4786.0386
$ws.Range("L126").Value = 16047.75
$ws.Range("M126").Value = -2316.0386
$ws.Range("N126").Value = -20987.75
$ws.Range("H134").Value = 6013.591
$ws.Range("I134").Value = 4893
$ws.Range("K134").Value = 14679
$ws.Range("M134").Value = -12144

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 467
$ws.Range("I8").Value = 467
$ws.Range("K8").Value = 1401
$ws.Range("M8").Value = -1262
$ws.Range("H11").Value = 8626.8125
$ws.Range("I11").Value = 9192.267
$ws.Range("J11").Value = 145
$ws.Range("K11").Value = 27576.801
$ws.Range("L11").Value = 435
$ws.Range("M11").Value = -27436.801
$ws.Range("N11").Value = -715
$ws.Range("H97").Value = 3174.6155
$ws.Range("I97").Value = 2046.375
$ws.Range("J97").Value = 4979.8
$ws.Range("K97").Value = 6139.125
$ws.Range("L97").Value = 14939.4
$ws.Range("M97").Value = -5643.125
$ws.Range("N97").Value = -15931.4
$ws.Range("H131").Value = 1517302.6
$ws.Range("J131").Value = 2275679.2
$ws.Range("L131").Value = 6827037.600000001
$ws.Range("N131").Value = -6837117.600000001
$ws.Range("H140").Value = 4430.5713
$ws.Range("I140").Value = 11316.7
$ws.Range("K140").Value = 33950.10000000001
$ws.Range("M140").Value = -28770.10000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1819.7778
$ws.Range("I2").Value = 1050.4667
$ws.Range("K2").Value = 1050.4667
$ws.Range("M2").Value = -937.4666999999999
$ws.Range("H70").Value = 5096.8887
$ws.Range("I70").Value = 4975.2
$ws.Range("K70").Value = 4975.2
$ws.Range("M70").Value = -4705.2
$ws.Range("H73").Value = 5096.8887
$ws.Range("I73").Value = 4975.2
$ws.Range("K73").Value = 4975.2
$ws.Range("M73").Value = -4039.2
$ws.Range("H113").Value = 3985.6316
$ws.Range("I113").Value = 2463.4
$ws.Range("K113").Value = 2463.4
$ws.Range("M113").Value = -293.4000000000001
$ws.Range("H126").Value = 7006.3
$ws.Range("J126").Value = 7310.4
$ws.Range("L126").Value = 21931.2
$ws.Range("N126").Value = -26871.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 73085.234
$ws.Range("I7").Value = 106954.82
$ws.Range("K7").Value = 106954.82
$ws.Range("M7").Value = -106842.82
$ws.Range("H22").Value = 1654.625
$ws.Range("I22").Value = 1068.7273
$ws.Range("K22").Value = 1068.7273
$ws.Range("M22").Value = -773.7273
$ws.Range("H27").Value = 1654.625
$ws.Range("I27").Value = 1068.7273
$ws.Range("K27").Value = 1068.7273
$ws.Range("M27").Value = -961.7273
$ws.Range("H46").Value = 4172.1665
$ws.Range("I46").Value = 2912.1667
$ws.Range("K46").Value = 2912.1667
$ws.Range("M46").Value = -2724.1667
$ws.Range("H55").Value = 1036.1428
$ws.Range("I55").Value = 689.5
$ws.Range("K55").Value = 689.5
$ws.Range("M55").Value = -516.5
$ws.Range("H126").Value = 73085.234
$ws.Range("I126").Value = 106954.82
$ws.Range("K126").Value = 320864.46
$ws.Range("M126").Value = -318394.46

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7521.35
$ws.Range("I62").Value = 5998.091
$ws.Range("J62").Value = 9383.111000000001
$ws.Range("K62").Value = 5998.091
$ws.Range("L62").Value = 9383.111000000001
$ws.Range("M62").Value = -5374.091
$ws.Range("N62").Value = -10631.111
$ws.Range("H64").Value = 157909.38
$ws.Range("J64").Value = 157909.38
$ws.Range("L64").Value = 157909.38
$ws.Range("N64").Value = -158405.38
$ws.Range("H65").Value = 7521.35
$ws.Range("I65").Value = 5998.091
$ws.Range("J65").Value = 9383.111000000001
$ws.Range("K65").Value = 29990.455
$ws.Range("L65").Value = 46915.55500000001
$ws.Range("M65").Value = -26870.455
$ws.Range("N65").Value = -53155.55500000001
$ws.Range("H67").Value = 157909.38
$ws.Range("J67").Value = 157909.38
$ws.Range("L67").Value = 157909.38
$ws.Range("N67").Value = -159625.38
$ws.Range("H122").Value = 4057
$ws.Range("I122").Value = 3626.923
$ws.Range("K122").Value = 10880.769
$ws.Range("M122").Value = -8430.769
